# Header modification for style ref
#
# This script applies a set of structural edits to the document:
#  1. Adds the new "Table header" paragraph style and its linked
#     "Table header Char" character style to the style sheet.
#  2. Removes the stray "_GoBack" bookmark that was sitting alone in an
#     empty Figure-caption paragraph.
#  3. Converts the two SEQ-Figure "fldSimple" fields into full complex
#     fields (begin / instrText / separate / result / end).
#  4. Splits the "Table header 1" cell text into three runs so that the
#     word "header" can carry the new TableheaderChar character style,
#     and re-adds a "_GoBack" bookmark at the end of that paragraph.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) New styles: paragraph style "Tableheader" + character style
#    "TableheaderChar", linked to one another.
# ---------------------------------------------------------------------
$paraStyle = $d.Styles.Add("Table header", 1)
$paraStyle.BaseStyle = "BodyText"

$charStyle = $d.Styles.Add("Table header Char", 2)
$charStyle.BaseStyle = "BodyTextChar"

$paraStyle.LinkStyle = "TableheaderChar"
$charStyle.LinkStyle = "Tableheader"
$paraStyle.QuickStyle = $true
$charStyle.Font.Name = "Arial Narrow"

# ---------------------------------------------------------------------
# 2) Remove the orphan "_GoBack" bookmark that lives alone inside an
#    empty Figure-caption paragraph (just before the bibliography).
# ---------------------------------------------------------------------
$goBackRange = $d.Bookmarks("_GoBack").Range
$goBackPara = $goBackRange.Paragraphs(1)
$emptyFigCaptionXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Figurecaption`"/></w:pPr></w:p>"
$goBackPara.Range.InsertXML($emptyFigCaptionXml) | Out-Null

# ---------------------------------------------------------------------
# 3) Convert the two SEQ Figure fldSimple fields to complex fields.
# ---------------------------------------------------------------------

# 3a. Figure caption paragraph (has the _Ref151977236 bookmark).
$figRange = $d.Content
$figRange.Find.ClearFormatting()
$figFound = $figRange.Find.Execute("Figure *. -- This is a figure caption*", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($figFound) {
    $figPara = $figRange.Paragraphs(1)
    $figXml = "<w:p $wNs>" + `
        "<w:pPr><w:pStyle w:val=`"Figurecaption`"/></w:pPr>" + `
        "<w:bookmarkStart w:id=`"100`" w:name=`"_Ref151977236`"/>" + `
        "<w:r><w:t xml:space=`"preserve`">Figure </w:t></w:r>" + `
        "<w:r><w:fldChar w:fldCharType=`"begin`"/></w:r>" + `
        "<w:r><w:instrText xml:space=`"preserve`"> SEQ Figure \* ARABIC </w:instrText></w:r>" + `
        "<w:r><w:fldChar w:fldCharType=`"separate`"/></w:r>" + `
        "<w:r><w:t>1</w:t></w:r>" + `
        "<w:r><w:fldChar w:fldCharType=`"end`"/></w:r>" + `
        "<w:bookmarkEnd w:id=`"100`"/>" + `
        "<w:r><w:t>. -- This is a figure caption. There is a weird reverse indent thing that only happens in tech memos</w:t></w:r>" + `
        "<w:r><w:t xml:space=`"preserve`"> and DPRs</w:t></w:r>" + `
        "<w:r><w:t>.</w:t></w:r>" + `
        "</w:p>"
    $figPara.Range.InsertXML($figXml) | Out-Null
}

# 3b. Image caption paragraph.
$imgRange = $d.Content
$imgRange.Find.ClearFormatting()
$imgFound = $imgRange.Find.Execute("Image *. -- This is an*", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($imgFound) {
    $imgPara = $imgRange.Paragraphs(1)
    $imgXml = "<w:p $wNs>" + `
        "<w:pPr><w:pStyle w:val=`"Imagecaption0`"/></w:pPr>" + `
        "<w:r><w:t>Image</w:t></w:r>" + `
        "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" + `
        "<w:r><w:fldChar w:fldCharType=`"begin`"/></w:r>" + `
        "<w:r><w:instrText xml:space=`"preserve`"> SEQ Figure \* ARABIC </w:instrText></w:r>" + `
        "<w:r><w:fldChar w:fldCharType=`"separate`"/></w:r>" + `
        "<w:r><w:t>1</w:t></w:r>" + `
        "<w:r><w:fldChar w:fldCharType=`"end`"/></w:r>" + `
        "<w:r><w:t xml:space=`"preserve`">. -- This is an " + [char]0x201C + "image caption." + [char]0x201D + " I" + [char]0x2019 + "m not sure if it matters but I" + [char]0x2019 + "m including it here. </w:t></w:r>" + `
        "<w:r><w:t>Here is some extra text that will show the indentation thing.</w:t></w:r>" + `
        "</w:p>"
    $imgPara.Range.InsertXML($imgXml) | Out-Null
}

# ---------------------------------------------------------------------
# 4) Split "Table header 1" into three runs (plain / styled / plain)
#    and re-add the "_GoBack" bookmark at the end of the paragraph.
# ---------------------------------------------------------------------
$tblRange = $d.Content
$tblRange.Find.ClearFormatting()
$tblFound = $tblRange.Find.Execute("Table header 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($tblFound) {
    $tblPara = $tblRange.Paragraphs(1)
    $tblXml = "<w:p $wNs>" + `
        "<w:pPr><w:pStyle w:val=`"BodyText`"/></w:pPr>" + `
        "<w:r><w:t xml:space=`"preserve`">Table </w:t></w:r>" + `
        "<w:r><w:rPr><w:rStyle w:val=`"TableheaderChar`"/></w:rPr><w:t>header</w:t></w:r>" + `
        "<w:r><w:t xml:space=`"preserve`"> 1</w:t></w:r>" + `
        "<w:bookmarkStart w:id=`"101`" w:name=`"_GoBack`"/>" + `
        "<w:bookmarkEnd w:id=`"101`"/>" + `
        "</w:p>"
    $tblPara.Range.InsertXML($tblXml) | Out-Null
}
